$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("2025-11-17")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2025-11-24"

# Header row (values + formatting copied from the previous week sheet)
$ws.Cells.Item(1,1).Value = 'rank'
$ws.Cells.Item(1,2).Value = 'title'
$ws.Cells.Item(1,3).Value = 'author'
$ws.Cells.Item(1,4).Value = 'latest_episode'
$src.Range("A1:D1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)

# Data rows (rank, title, author, latest_episode)
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = '新米オッサン冒険者、最強パーティに死ぬほど鍛えられて無敵になる'
$ws.Cells.Item(2,3).Value = '漫画：荻野ケン 原作：岸馬きらく キャラクター原案：Tea'
$ws.Cells.Item(2,4).Value = '第73話 前編'
$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = '勇者パーティーをクビになったので故郷に帰ったら、メンバー全員がついてきたんだが'
$ws.Cells.Item(3,3).Value = '絶叫あいす。(漫画) 木の芽(原作) 希(キャラクター原案)'
$ws.Cells.Item(3,4).Value = 'コミックス第1巻 発売告知記事'
$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = 'いとこのこ'
$ws.Cells.Item(4,3).Value = 'いぬちく(著者)'
$ws.Cells.Item(4,4).Value = '第42話'
$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = '時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―'
$ws.Cells.Item(5,3).Value = '光永康則'
$ws.Cells.Item(5,4).Value = '第７１話『扇山停止』③'
$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = '生徒会にも穴はある！'
$ws.Cells.Item(6,3).Value = 'むちまろ'
$ws.Cells.Item(6,4).Value = '第140話	蚊がせめてきたぞっ!!'
$ws.Cells.Item(7,1).Value = 6
$ws.Cells.Item(7,2).Value = '地元のいじめっ子達に仕返ししようとしたら、別の戦いが始まった。'
$ws.Cells.Item(7,3).Value = 'マツモトケンゴ'
$ws.Cells.Item(7,4).Value = '第６９話　爺ちゃんとの戦いが始まった（１）'
$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = '元・世界１位のサブキャラ育成日記 ～廃プレイヤー、異世界を攻略中！～'
$ws.Cells.Item(8,3).Value = '沢村治太郎(原作) 前田理想(漫画) まろ(キャラクター原案)'
$ws.Cells.Item(8,4).Value = '第80話その2'
$ws.Cells.Item(9,1).Value = 8
$ws.Cells.Item(9,2).Value = '実は俺、最強でした？'
$ws.Cells.Item(9,3).Value = '原作：澄守 彩 漫画：高橋 愛'
$ws.Cells.Item(9,4).Value = '第132話　画伯の願い'
$ws.Cells.Item(10,1).Value = 9
$ws.Cells.Item(10,2).Value = '怠惰な悪辱貴族に転生した俺、シナリオをぶっ壊したら規格外の魔力で最凶になった'
$ws.Cells.Item(10,3).Value = '菊池快晴(原作) 小田童馬(作画) 桑島黎音(キャラクター原案)'
$ws.Cells.Item(10,4).Value = '第16話'
$ws.Cells.Item(11,1).Value = 10
$ws.Cells.Item(11,2).Value = '女友達は頼めば意外とヤらせてくれる'
$ws.Cells.Item(11,3).Value = 'ろくろ(漫画) 鏡遊(原作)'
$ws.Cells.Item(11,4).Value = '第26話②'
$ws.Cells.Item(12,1).Value = 11
$ws.Cells.Item(12,2).Value = '異世界魔王と召喚少女の奴隷魔術'
$ws.Cells.Item(12,3).Value = '原作：むらさきゆきや 漫画：福田直叶 キャラクター原案：鶴崎貴大'
$ws.Cells.Item(12,4).Value = '第130話　変身してみる（後編）'
$ws.Cells.Item(13,1).Value = 12
$ws.Cells.Item(13,2).Value = '勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～'
$ws.Cells.Item(13,3).Value = '漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり'
$ws.Cells.Item(13,4).Value = '第５４話　勇者を護る器用貧乏（４）'
$ws.Cells.Item(14,1).Value = 13
$ws.Cells.Item(14,2).Value = '異世界のんびり農家'
$ws.Cells.Item(14,3).Value = '剣康之(作画) 内藤騎之介(原作) やすも(キャラクター原案)'
$ws.Cells.Item(14,4).Value = '第311話'
$ws.Cells.Item(15,1).Value = 14
$ws.Cells.Item(15,2).Value = '氷結令嬢さまをフォローしたら、メチャメチャ溺愛されてしまった件@comic'
$ws.Cells.Item(15,3).Value = '漫画：ハレノチアメ 原作：愛坂タカト キャラクター原案：Bcoca'
$ws.Cells.Item(15,4).Value = '第10話'
$ws.Cells.Item(16,1).Value = 15
$ws.Cells.Item(16,2).Value = 'バキ外伝 烈海王は異世界転生しても一向にかまわんッッ'
$ws.Cells.Item(16,3).Value = '板垣恵介 猪原賽 陸井栄史'
$ws.Cells.Item(16,4).Value = '第84話　慈愛'
$ws.Cells.Item(17,1).Value = 16
$ws.Cells.Item(17,2).Value = '【パクパクですわ】追放されたお嬢様の『モンスターを食べるほど強くなる』スキルは、１食で１レベルアップする前代未聞の最強スキルでした。３日で人類最強になりましたわ～！'
$ws.Cells.Item(17,3).Value = '島知宏 音速炒飯 有都あらゆる'
$ws.Cells.Item(17,4).Value = '第２５食　赤スライムのシャーベット、パクパクですわ！（３）'
$ws.Cells.Item(18,1).Value = 17
$ws.Cells.Item(18,2).Value = 'リビルドワールド'
$ws.Cells.Item(18,3).Value = '綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)'
$ws.Cells.Item(18,4).Value = '第75話③'
$ws.Cells.Item(19,1).Value = 18
$ws.Cells.Item(19,2).Value = 'ぽんドロイド！ はまさん'
$ws.Cells.Item(19,3).Value = 'はれやまはれぞう(著者)'
$ws.Cells.Item(19,4).Value = '第12話'
$ws.Cells.Item(20,1).Value = 19
$ws.Cells.Item(20,2).Value = 'アザミヤコを好きになる'
$ws.Cells.Item(20,3).Value = 'ユニティコング(原作) ツノニガウ(作画)'
$ws.Cells.Item(20,4).Value = '第11話後編'
$ws.Cells.Item(21,1).Value = 20
$ws.Cells.Item(21,2).Value = '美人女上司滝沢さん'
$ws.Cells.Item(21,3).Value = 'やんBARU(著者)'
$ws.Cells.Item(21,4).Value = '第208話'
$ws.Cells.Item(22,1).Value = 21
$ws.Cells.Item(22,2).Value = '貞操逆転世界で頼めばヤれると噂の俺'
$ws.Cells.Item(22,3).Value = '澄田佑貴(漫画) aaa168（スリーエー）(原作)'
$ws.Cells.Item(22,4).Value = '第1話'
$ws.Cells.Item(23,1).Value = 22
$ws.Cells.Item(23,2).Value = '世界最強の魔女、始めました 〜私だけ『攻略サイト』を見れる世界で自由に生きます〜'
$ws.Cells.Item(23,3).Value = '戸賀 環 坂木持丸 riritto'
$ws.Cells.Item(23,4).Value = '第56話②　ペットを飼ってみた'
$ws.Cells.Item(24,1).Value = 23
$ws.Cells.Item(24,2).Value = 'よくわからないけれど異世界に転生していたようです'
$ws.Cells.Item(24,3).Value = '内々けやき あし カオミン'
$ws.Cells.Item(24,4).Value = '第143話 よくわからないけれど人をダメにするみたいです（１）'
$ws.Cells.Item(25,1).Value = 24
$ws.Cells.Item(25,2).Value = '聖者無双'
$ws.Cells.Item(25,3).Value = '漫画：秋風緋色 原作：ブロッコリーライオン キャラクター原案：sime'
$ws.Cells.Item(25,4).Value = '第94話　戦乱のドワーフ王国・奴隷の扱い（後半）'
$ws.Cells.Item(26,1).Value = 25
$ws.Cells.Item(26,2).Value = '異世界メイドの三ツ星グルメ ～現代ごはん作ったら王宮で大バズリしました～'
$ws.Cells.Item(26,3).Value = 'モリタ Ｕ４ nima'
$ws.Cells.Item(26,4).Value = '第14話（３）　春とおぼっちゃまとピクニックランチ（３）'
$ws.Cells.Item(27,1).Value = 26
$ws.Cells.Item(27,2).Value = '落ちこぼれだった兄が実は最強 ～史上最強の勇者は転生し、学園で無自覚に無双する～'
$ws.Cells.Item(27,3).Value = '村上よしゆき 茨木野 あるてら'
$ws.Cells.Item(27,4).Value = '第４３話　勇者、合体した六邪神将を撃破し、めでたしめでたし（２）'
$ws.Cells.Item(28,1).Value = 27
$ws.Cells.Item(28,2).Value = '配信に致命的に向いていない女の子が迷宮で黙々と人助けする配信'
$ws.Cells.Item(28,3).Value = '下田将也(漫画) 佐藤悪糖(原作) 福きつね(キャラクター原案)'
$ws.Cells.Item(28,4).Value = '第4話後編'
$ws.Cells.Item(29,1).Value = 28
$ws.Cells.Item(29,2).Value = '小林さんちのメイドラゴン'
$ws.Cells.Item(29,3).Value = 'クール教信者'
$ws.Cells.Item(29,4).Value = '第154話'
$ws.Cells.Item(30,1).Value = 29
$ws.Cells.Item(30,2).Value = 'バキ外伝　ガイアとシコルスキー　～ときどきノムラ 二人だけど三人暮らし～'
$ws.Cells.Item(30,3).Value = '板垣恵介 林たかあき'
$ws.Cells.Item(30,4).Value = '第58話 帰還'
$ws.Cells.Item(31,1).Value = 30
$ws.Cells.Item(31,2).Value = 'くらいあの子としたいこと'
$ws.Cells.Item(31,3).Value = '碇マナツ(著者)'
$ws.Cells.Item(31,4).Value = '特別編㉒'
$ws.Cells.Item(32,1).Value = 31
$ws.Cells.Item(32,2).Value = 'ダンジョンの幼なじみ'
$ws.Cells.Item(32,3).Value = '久真やすひさ(著者)'
$ws.Cells.Item(32,4).Value = '【７巻発売＆1000万PV突破記念！】 ダンジョンの幼なじみ第２回人気投票'
$ws.Cells.Item(33,1).Value = 32
$ws.Cells.Item(33,2).Value = '姫様“拷問”の時間です'
$ws.Cells.Item(33,3).Value = '原作:春原ロビンソン　漫画:ひらけい'
$ws.Cells.Item(33,4).Value = '拷問158'
$ws.Cells.Item(34,1).Value = 33
$ws.Cells.Item(34,2).Value = '追放されたチート付与魔術師は 気ままなセカンドライフを謳歌する。'
$ws.Cells.Item(34,3).Value = '六志麻あさ 業務用餅 kisui'
$ws.Cells.Item(34,4).Value = '第７５話'
$ws.Cells.Item(35,1).Value = 34
$ws.Cells.Item(35,2).Value = '十年目、帰還を諦めた転移者はいまさら主人公になる'
$ws.Cells.Item(35,3).Value = '原作：氷純（「十年目、帰還を諦めた転移者はいまさら主人公になる」MFブックス刊） 漫画：しゅーかま キャラクター原案：あんべよしろう'
$ws.Cells.Item(35,4).Value = '第20話④'
$ws.Cells.Item(36,1).Value = 35
$ws.Cells.Item(36,2).Value = 'お気楽領主の楽しい領地防衛 ～生産系魔術で名もなき村を最強の城塞都市に～'
$ws.Cells.Item(36,3).Value = '青色まろ（漫画） 赤池宗（原作） 転（原作イラスト）'
$ws.Cells.Item(36,4).Value = '第35話　侵略者'
$ws.Cells.Item(37,1).Value = 36
$ws.Cells.Item(37,2).Value = '理想のヒモ生活'
$ws.Cells.Item(37,3).Value = '日月ネコ(漫画) 渡辺恒彦（ヒーロー文庫／イマジカインフォス）(原作) 文倉十(キャラクター原案)'
$ws.Cells.Item(37,4).Value = '第89話　その3'
$ws.Cells.Item(38,1).Value = 37
$ws.Cells.Item(38,2).Value = 'ふかふかダンジョン攻略記～俺の異世界転生冒険譚～'
$ws.Cells.Item(38,3).Value = 'KAKERU'
$ws.Cells.Item(38,4).Value = '第70話「あうと！ せーふ！ よよいのよい！（もうどうにでもなぁれ♡）」（後半)'
$ws.Cells.Item(39,1).Value = 38
$ws.Cells.Item(39,2).Value = '最弱貴族に転生したので悪役たちを集めてみた'
$ws.Cells.Item(39,3).Value = '空野進 sorani ファルまろ'
$ws.Cells.Item(39,4).Value = '第14話　最弱貴族、悪役令嬢を脱がす（２）'
$ws.Cells.Item(40,1).Value = 39
$ws.Cells.Item(40,2).Value = '最果てのパラディン'
$ws.Cells.Item(40,3).Value = '奥橋睦（漫画） 柳野かなた（原作） 輪くすさが（キャラクター原案）'
$ws.Cells.Item(40,4).Value = '第68話　無敵の巨人Ⅰ'
$ws.Cells.Item(41,1).Value = 40
$ws.Cells.Item(41,2).Value = 'Lv２からチートだった元勇者候補のまったり異世界ライフ'
$ws.Cells.Item(41,3).Value = '糸町秋音（漫画） 鬼ノ城ミヤ（原作） 片桐（キャラクター原案）'
$ws.Cells.Item(41,4).Value = '第63話　居場所'
$ws.Cells.Item(42,1).Value = 41
$ws.Cells.Item(42,2).Value = '黄金の経験値'
$ws.Cells.Item(42,3).Value = '原純(原作) 霜月汐(作画) fixro2n(キャラクター原案)'
$ws.Cells.Item(42,4).Value = '第20話（前編）'
$ws.Cells.Item(43,1).Value = 42
$ws.Cells.Item(43,2).Value = '婚約者に裏切られた錬金術師は、独立して『ざまぁ』します　コミック版'
$ws.Cells.Item(43,3).Value = '漫画/すたひろ 原作/Y.A'
$ws.Cells.Item(43,4).Value = 'chapter74【39話①】'
$ws.Cells.Item(44,1).Value = 43
$ws.Cells.Item(44,2).Value = '俺以外誰も採取できない素材なのに「素材採取率が低い」とパワハラする幼馴染錬金術師と絶縁した専属魔導士、辺境の町でスローライフを送りたい。'
$ws.Cells.Item(44,3).Value = '狐御前(原作) 西岡知三(作画) ＮＯＣＯ(キャラクター原案)'
$ws.Cells.Item(44,4).Value = '第27話-2'
$ws.Cells.Item(45,1).Value = 44
$ws.Cells.Item(45,2).Value = '塔の管理をしてみよう'
$ws.Cells.Item(45,3).Value = '盧恩＆雪笠(Friendly Land)(著者) 早秋(原作) 雨神(キャラクター原案)'
$ws.Cells.Item(45,4).Value = '第95話前編'
$ws.Cells.Item(46,1).Value = 45
$ws.Cells.Item(46,2).Value = '10年ぶりに再会したクソガキは清純美少女JKに成長していた'
$ws.Cells.Item(46,3).Value = '緑青黒羽（漫画） 館西夕木（原作） ひげ猫（キャラクター原案）'
$ws.Cells.Item(46,4).Value = '第7話　ショッピングデート（後編）'
$ws.Cells.Item(47,1).Value = 46
$ws.Cells.Item(47,2).Value = '男子高校生だけどギャルにTSしました'
$ws.Cells.Item(47,3).Value = '太陽まりい(著者)'
$ws.Cells.Item(47,4).Value = '第21話前編'
$ws.Cells.Item(48,1).Value = 47
$ws.Cells.Item(48,2).Value = '異世界のすみっこで快適ものづくり生活 ～女神さまのくれた工房はちょっとやりすぎ性能だった～'
$ws.Cells.Item(48,3).Value = '西山アラタ(漫画) 長田信織(原作) 東上文(キャラクター原案)'
$ws.Cells.Item(48,4).Value = 'EP.21②'
$ws.Cells.Item(49,1).Value = 48
$ws.Cells.Item(49,2).Value = '理想の彼女'
$ws.Cells.Item(49,3).Value = 'もりまりも(著者)'
$ws.Cells.Item(49,4).Value = '第29話'
$ws.Cells.Item(50,1).Value = 49
$ws.Cells.Item(50,2).Value = '魔法少女リリカルなのは EXCEEDS'
$ws.Cells.Item(50,3).Value = '都築真紀 川上修一'
$ws.Cells.Item(50,4).Value = '第７話①'
$ws.Cells.Item(51,1).Value = 50
$ws.Cells.Item(51,2).Value = '魔石グルメ　魔物の力を食べたオレは最強！'
$ws.Cells.Item(51,3).Value = '菅原健二(作画) 結城涼(原作) 成瀬ちさと(キャラクター原案)'
$ws.Cells.Item(51,4).Value = '第69話後半'

# Restore original active sheet/selection state
$wb.Worksheets.Item("Sheet1").Activate()
